# Append the newest (7th October 2024) rows to the 60000-80000 analysis sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT even when it looks like a pure
# number (e.g. CVR numbers), without leaving any new/unused cell style
# behind. We do this by building the text via a formula (which always
# yields a text-typed result), copying it, and pasting *values only* into
# the destination - Excel keeps the text type without needing a text
# number format or a quote-prefix on the destination cell.
function Set-TextValue {
    param($cell, [string]$text)
    $scratch = $ws.Range("Z100")
    $escaped = $text.Replace("""", """""")
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Date cells in column E share the same custom date number format as the
# existing rows - grab it once so new date cells reuse the same style.
$dateFormat = $ws.Range("E2").NumberFormat

$rows = @(
    @{ Row=30; A="45445216"; B=2024; C=76625;    D="Visma Løn";  E=45301; G=$null;             H="2024Q1" },
    @{ Row=31; A="35651950"; B=2024; C=66672;    D="Visma Time"; E=45407; G=$null;             H="2024Q2" },
    @{ Row=32; A="30823699"; B=2024; C=78804;    D="EasyCruit";  E=45469; G="Vil ikke oplyse"; H="2024Q2" },
    @{ Row=33; A="29616647"; B=2024; C=60060;    D="Visma Løn";  E=45498; G=$null;             H="2024Q3" },
    @{ Row=34; A="10154529"; B=2024; C=67236;    D="Visma Løn";  E=45513; G=$null;             H="2024Q3" },
    @{ Row=35; A="24997189"; B=2024; C=69000;    D="EasyCruit";  E=45551; G=$null;             H="2024Q3" },
    @{ Row=36; A="10658446"; B=2024; C=61860;    D="EasyCruit";  E=45561; G=$null;             H="2024Q3" }
)

foreach ($r in $rows) {
    $row = $r.Row

    Set-TextValue $ws.Cells.Item($row, 1) $r.A   # A: CVR (numeric-looking -> must stay text)
    $ws.Cells.Item($row, 2).Value = $r.B          # B: Year
    $ws.Cells.Item($row, 3).Value = $r.C          # C: Beløb 12 mdr. (TCV)
    $ws.Cells.Item($row, 4).Value = $r.D          # D: Løsning

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = $dateFormat
    $eCell.Value = $r.E                            # E: Opsagt dato:

    if ($r.G) {
        $ws.Cells.Item($row, 7).Value = $r.G       # G: Ny leverandør
    }

    $ws.Cells.Item($row, 8).Value = $r.H           # H: Quarter
    $ws.Cells.Item($row, 9).Value = "60000-80000"  # I: TCV_range
}
